$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows (Ca8/heme and Ca9/solid) right after the existing
# "Ca7 / cancer_status" row (row 8), before the old row 9 (D1/age). ---
$ws.Rows("9:10").Insert()

$ws.Range("A9").Value = "Ca8"
$ws.Range("B9").Value = "heme"
$ws.Range("C9").Value = "Cancer"
$ws.Range("D9").Value = "Hematologic malignancy indicator"

$ws.Range("A10").Value = "Ca9"
$ws.Range("B10").Value = "solid"
$ws.Range("C10").Value = "Cancer"
$ws.Range("D10").Value = "Solid tumor indicator"

# Grow the existing table (Table1) so it covers the two newly inserted rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E68"))

# --- Append a new row at the bottom of the table (X2/imwg) ---
$lo.ListRows.Add()

$ws.Range("A69").Value = "X2"
$ws.Range("B69").Value = "imwg"
$ws.Range("C69").Value = "Other"
$ws.Range("D69").Value = "Modified IMWG frailty index"

# Restore the selection to reflect where the user ended up after the edits.
$ws.Range("A71").Select()
